$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 5 (Quarta / 17h40 / 19h), pushing nothing else around
# since row 12 stays at row 12 in the target (no shift - row 5 was empty before).
$ws.Range("A5").Value = "Quarta"
$ws.Range("B5").Value = "17h40"
$ws.Range("C5").Value = "19h"

# Update the active selection to C6, matching the saved workbook view.
$ws.Range("C6").Select()
